# Improved reward function: lower the "intended direction" probability
# on the "normal" movement-probabilities sheet from 0.9 to 0.8.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("normal")

$ws.Range("B2").Value = 0.8
$ws.Range("C3").Value = 0.8
$ws.Range("D4").Value = 0.8
$ws.Range("E5").Value = 0.8
$ws.Range("F6").Value = 0.8
$ws.Range("G7").Value = 0.8
$ws.Range("H8").Value = 0.8
$ws.Range("I9").Value = 0.8

$ws.Range("I9").Select()
